$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (rows 2-51) to Text format before writing the new
# values. Several values look like numbers/dates to Excel (e.g. "1.000",
# "1.001", "0.9996") and would otherwise get silently reinterpreted/
# reformatted as numeric values, losing significant trailing zeros and the
# original string type. Restoring the Normal style afterwards keeps the
# cells free of an explicit style index (matching the source workbook,
# where these cells carry no "s" attribute) while preserving the text type.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '25.811.00'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '1.739.53'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '227.58'
$ws.Range('E5').Value = '  -3.84%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '0.5159'
$ws.Range('E7').Value = '  +2.09%  '
$ws.Range('D8').Value = '0.2724'
$ws.Range('E8').Value = '  +2.86%  '
$ws.Range('D9').Value = '38.76'
$ws.Range('E9').Value = '  -5.80%  '
$ws.Range('D10').Value = '0.06085'
$ws.Range('E10').Value = '  -2.03%  '
$ws.Range('D11').Value = '1.738.54'
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('D12').Value = '0.07003'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').Value = '15.17'
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('D14').Value = '0.6285'
$ws.Range('E14').Value = '  +4.69%  '
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').Value = '76.23'
$ws.Range('E16').Value = '  -1.20%  '
$ws.Range('D17').Value = '0.9997'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '0.9996'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = '25.829.68'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('D21').Value = '0.000006602'
$ws.Range('E21').Value = '  -3.17%  '
$ws.Range('D22').Value = '1.960.35'
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').Value = '4.045'
$ws.Range('E23').Value = '  -0.77%  '
$ws.Range('D24').Value = '8.431'
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('E25').Value = '  -1.78%  '
$ws.Range('D26').Value = '136.58'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = '1.506'
$ws.Range('E27').Value = '  +4.19%  '
$ws.Range('D28').Value = '1.819'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').Value = '14.97'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').Value = '102.49'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').Value = '0.08315'
$ws.Range('E31').Value = '  +1.67%  '
$ws.Range('E32').Value = '  -1.41%  '
$ws.Range('D33').Value = '3.374'
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('D34').Value = '0.04407'
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('D35').Value = '2.605'
$ws.Range('E35').Value = '  -1.81%  '
$ws.Range('D36').Value = '0.9698'
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('D37').Value = '0.5949'
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('D38').Value = '2.681'
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('D39').Value = '0.01558'
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('D40').Value = '1.932'
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('D41').Value = '0.9989'
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').Value = '101.72'
$ws.Range('E42').Value = '  -1.51%  '
$ws.Range('D43').Value = '0.3795'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('D44').Value = '0.7245'
$ws.Range('E44').Value = '  -2.14%  '
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('D46').Value = '0.05482'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '6.183'
$ws.Range('E47').Value = '  +4.06%  '
$ws.Range('D48').Value = '0.1097'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').Value = '29.82'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').Value = '51.86'
$ws.Range('E50').Value = '  -0.53%  '
$ws.Range('E51').Value = '  +0.24%  '

$priceRange.Style = "Normal"
